$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet, so it lands after "Lookup".
$lastIndex = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$ws.Name = "LookupBlanks"

# Header row (row 4)
$ws.Range("D4").Value = "VLOOKUP"
$ws.Range("E4").Value = "MATCH"

# Column A - lookup table values (row 8 intentionally left blank to exercise
# VLOOKUP/MATCH/LOOKUP behaviour against blank cells in the lookup range).
$ws.Range("A5").Value = 1
$ws.Range("A6").Value = 2
$ws.Range("A7").Value = 3
$ws.Range("A9").Value = 5
$ws.Range("A10").Value = 6
$ws.Range("A11").Value = 7
$ws.Range("A12").Value = 8

# Column C - lookup values
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 2
$ws.Range("C8").Value = 6
$ws.Range("C9").Value = 7
$ws.Range("C10").Value = 8
$ws.Range("C11").Value = 9
$ws.Range("C12").Value = "'Hi"

# F4 header written after C12 so the shared-string table order matches the
# authored workbook (Hi before LOOKUP).
$ws.Range("F4").Value = "LOOKUP"

# Formulas (rows 5-12), shared down each column.
$ws.Range("D5:D12").Formula = "=VLOOKUP(C5,A`$1:A`$16,1)"
$ws.Range("E5:E12").Formula = "=MATCH(C5,A`$1:A`$16)"
$ws.Range("F5:F12").Formula = "=LOOKUP(C5,A`$1:A`$16)"

# Re-select the original "Lookup" sheet so it stays the active/visible tab.
$wb.Worksheets.Item(1).Activate()
